# Apply the Data1.xlsx update:
#  1) Refresh the POC test data values used in row 5 of Sheet1
#  2) Move the active selection to G7
#  3) Set the page to Portrait orientation (adds a pageSetup element)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the test-data row (row 5). A5/B5 share one value, C5/D5 share
#    another, E5/F5 share a third, and H5 gets the new "Project" value.
#    H5 is written first to match the authored edit order.
$ws.Range("H5").Value = "Project222"
$ws.Range("A5").Value = "POC5OCTo"
$ws.Range("B5").Value = "POC5OCTo"
$ws.Range("C5").Value = "FacilityPOC5OCTo"
$ws.Range("D5").Value = "FacilityPOC5OCTo"
$ws.Range("E5").Value = "PharmacyPOC5OCTo"
$ws.Range("F5").Value = "PharmacyPOC5OCTo"
# G5 (Cerner) is unchanged.

# 2) Move/save the selected cell as G7 (was H5).
$ws.Range("G7").Select()

# 3) Explicitly set Portrait orientation so a <pageSetup> entry is written.
$ws.PageSetup.Orientation = 1
